# Auto-generated edit script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format on the Price/Volume columns so the
# numeric-looking strings (e.g. "60.896.52", "0.994") are not
# auto-converted to numbers/dates by Excel's cell-value parser.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.896.52"
$ws.Range("E2").Value = "  -2.38%  "
$ws.Range("D3").Value = "2.420.11"
$ws.Range("E3").Value = "  -1.37%  "
$ws.Range("D4").Value = "0.994"
$ws.Range("E4").Value = "  -0.46%  "
$ws.Range("D5").Value = "570.52"
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").Value = "140.11"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").Value = "2.405.27"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("D11").Value = "0.159"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("E15").Value = "  -2.53%  "
$ws.Range("D16").Value = "2.852.16"
$ws.Range("E16").Value = "  -1.52%  "
$ws.Range("D17").Value = "60.800.89"
$ws.Range("E17").Value = "  -2.35%  "
$ws.Range("D18").Value = "2.401.05"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").Value = "7.47"
$ws.Range("E19").Value = "  +4.81%  "
$ws.Range("E20").Value = "  -2.17%  "
$ws.Range("D21").Value = "322.91"
$ws.Range("E21").Value = "  -1.68%  "
$ws.Range("D22").Value = "4.06"
$ws.Range("E22").Value = "  -1.53%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "1.87"
$ws.Range("E25").Value = "  -4.65%  "
$ws.Range("D26").Value = "64.90"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "579.62"
$ws.Range("D28").Value = "8.43"
$ws.Range("E28").Value = "  -10.06%  "
$ws.Range("D29").Value = "2.524.54"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "0.0₃0921"
$ws.Range("E30").Value = "  -4.10%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("E32").Value = "  -5.98%  "
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("D34").Value = "0.132"
$ws.Range("E34").Value = "  -3.16%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  -5.96%  "
$ws.Range("D37").Value = "1.40"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").Value = "150.10"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("D40").Value = "18.28"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Value = "5.14"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").Value = "1.68"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").Value = "41.13"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("E45").Value = "  -4.91%  "
$ws.Range("E46").Value = "  +14.88%  "
$ws.Range("D47").Value = "141.44"
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "3.52"
$ws.Range("E48").Value = "  -3.68%  "
$ws.Range("E49").Value = "  -2.86%  "
$ws.Range("D50").Value = "19.56"
$ws.Range("E50").Value = "  -1.20%  "
$ws.Range("E51").Value = "  -3.36%  "

# Restore the default (unstyled) cell style so formatting matches the original workbook.
$ws.Range("D2:E51").Style = "Normal"

Write-Host "Applied 78 cell updates to cryptos sheet"
